# Apply German translations to Menu-Languages.docx
$d = $word.ActiveDocument

function Replace-Text($find, $replace, $matchWholeWord) {
    $d.Content.Find.Execute($find, $matchWholeWord, $true, $false, $false, $false,
                             $true, 1, $false, $replace, 2)
}

Replace-Text "Menu:" "Menü:" $false
Replace-Text "Welcome" "Willkommen" $true
Replace-Text "What is SmartCash?" "Was ist SmartCash?" $false
Replace-Text "Brochure/Whitepaper" "Broschüre/Whitepaper" $false
Replace-Text "SmartHive Discussion" "SmartHive Diskussion" $false
Replace-Text "SmartHive Voting" "SmartHive Abstimmung" $false
Replace-Text "Publications Archive" "Publikations-Archiv" $false

Replace-Text "Businesses" "Unternehmen" $true
Replace-Text "Benefit for Merchants" "Vorteile für Händler" $false

# "SmartCard" occurs twice in the document; only the first (menu entry)
# should be translated, so replace just the first occurrence.
$d.Content.Find.Execute("SmartCard", $true, $true, $false, $false, $false,
                         $true, 1, $false, "Smartcard", 1)

Replace-Text "The Other Side" "Die andere Seite" $false
Replace-Text "Resources" "Ressourcen" $true
Replace-Text "Pool overview" "Pool Überblick" $false
Replace-Text "Services" "Dienstleistungen" $true
Replace-Text "Projects" "Projekte" $true

# "Exchanges Listing Guide" must be replaced before the standalone
# "Exchanges" entry, otherwise the standalone replace would also match
# the start of this longer phrase.
Replace-Text "Exchanges Listing Guide" "Börsen Listing Guide" $true
Replace-Text "Exchanges" "Börsen" $true

Replace-Text "Graphics" "Grafiken" $true
